# Applies the "PO Forecast" update:
#  1. Renames the header cells on the existing sheets.
#  2. Adds a new "PO Forecast" worksheet at the end of the workbook.
#  3. Fills it with the ds / PO_Forecast / yhat_lower / yhat_upper data,
#     reusing the header style and date-format style already used on the
#     other sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the "Requested quantity" headers on the existing sheets
# ---------------------------------------------------------------------------
$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# ---------------------------------------------------------------------------
# 2. Add the new "PO Forecast" worksheet after the last existing sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# ---------------------------------------------------------------------------
# Forecast data: ds, PO_Forecast, yhat_lower, yhat_upper
# ---------------------------------------------------------------------------
$poData = @(
    ,@(45235.99999999999, 54, 10.56122246714633, 95.5280123400338)
    ,@(45242.99999999999, 54, 16.01039474724748, 92.85033221222325)
    ,@(45270.99999999999, 56, 14.25682420722001, 92.77656284944391)
    ,@(45347.99999999999, 60, 20.76336610409046, 99.67628148464055)
    ,@(45382.99999999999, 62, 22.03209808588715, 102.9829143592871)
    ,@(45410.99999999999, 64, 21.80662240031462, 100.8829341929647)
    ,@(45417.99999999999, 64, 23.8503325393999, 104.9419388402644)
    ,@(45424.99999999999, 65, 24.8399517034643, 105.3100170577122)
    ,@(45438.99999999999, 65, 21.74499861470627, 104.5244714533646)
    ,@(45445.99999999999, 66, 24.35848410283221, 102.9273776785368)
    ,@(45452.99999999999, 66, 25.71196591466869, 104.5255047357416)
    ,@(45473.99999999999, 68, 27.38416644493781, 106.19636340965)
    ,@(45487.99999999999, 68, 27.2752734781396, 109.5914535639036)
    ,@(45501.99999999999, 69, 27.68748650577826, 109.8446004167103)
    ,@(45515.99999999999, 70, 28.77719306404674, 111.2516366617079)
    ,@(45522.99999999999, 70, 28.08427750026384, 111.9816253452518)
    ,@(45536.99999999999, 71, 29.70653448248467, 114.0244261237793)
    ,@(45543.99999999999, 72, 34.07833647493734, 112.7006524722513)
    ,@(45557.99999999999, 72, 32.56442390985968, 113.7941741632276)
    ,@(45571.99999999999, 73, 32.45490165529502, 112.8398516473206)
    ,@(45599.99999999999, 75, 33.23008030224883, 115.8645732856123)
    ,@(45606.99999999999, 75, 34.29490009621964, 115.0819619797318)
    ,@(45613.99999999999, 76, 33.27968250777142, 117.9355904139861)
    ,@(45620.99999999999, 76, 34.58020741737104, 119.8633617274456)
    ,@(45627.99999999999, 77, 37.49242536890365, 115.2667975440535)
    ,@(45634.99999999999, 77, 33.09512195067553, 119.5541231066535)
    ,@(45641.99999999999, 77, 37.56791917863551, 119.1095770395338)
    ,@(45648.99999999999, 78, 39.81153701808856, 117.0431588991999)
    ,@(45655.99999999999, 78, 36.45188856131641, 121.2754429153082)
    ,@(45662.99999999999, 79, 40.16423706158076, 118.7670735865966)
)

$lastRow = 1 + $poData.Length   # header row + data rows

# ---------------------------------------------------------------------------
# 3. Header row: write values, then stamp with the same bold/border/center
#    style already used for the headers on the other sheets.
# ---------------------------------------------------------------------------
$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 4. Data rows
# ---------------------------------------------------------------------------
$row = 2
foreach ($entry in $poData) {
    $wsForecast.Cells.Item($row, 1).Value = $entry[0]
    $wsForecast.Cells.Item($row, 2).Value = $entry[1]
    $wsForecast.Cells.Item($row, 3).Value = $entry[2]
    $wsForecast.Cells.Item($row, 4).Value = $entry[3]
    $row++
}

# Reuse the existing date-format style for the "ds" column, same as the
# "Order Week" / "Order Month" columns on the other sheets.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A$lastRow").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0
